$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update row 10 (Objetivos) B/C text with the new description
$ws.Range("B10:C10").Value = 'Apresentar a teoria clássica do controle automático de processos para análise e projeto 
de sistemas de controle feedback.'

# 2) Insert a brand new row at position 13 (shifts old rows 13-25 down to 14-26)
$ws.Rows.Item(13).Insert()

# 3) Populate the new row 13 (B/C only; A stays blank) with the lecturer name
$ws.Range("B13:C13").Value = '8643537 - Fabio Rodolfo Miguel Batista'
# Fix up B13 style (engine mis-assigns col-B style on brand new cells); copy correct format from B10
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) Update text in the shifted rows that also changed content
$ws.Range("B14:C14").Value = 'Introdução ao controle automático de processos. Modelagem de processos. Análise do 
comportamento dinâmico de processos. Análise e projeto de sistemas de controle feedback.'
$ws.Range("B16:C16").Value = '1- INTRODUÇÃO AO CONTROLE AUTOMÁTICO DE PROCESSOS. - Introdução à teoria de controle de processos. - Malhas aberta e fechada. - Diagrama de blocos.
2- MODELAGEM DE PROCESSOS. - Considerações de modelagem matemática para fins de controle de processos. - Transformada de Laplace. - Função de transferência e modelos entrada-saída.
3- ANÁLISE DO COMPORTAMENTO DINÂMICO DE PROCESSOS. - Sistemas de 1ª ordem. - Sistemas de 2ª ordem e ordem superior.
4- ANÁLISE E PROJETO DE SISTEMAS DE CONTROLE. - Comportamento dinâmico de processos controlados por sistemas de controle feedback. - Análise de estabilidade. Sintonia de controladores de realimentação.
5- RESPOSTA EM FREQUÊNCIA. - Análise de sistemas de controle através de resposta em frequencia.'
$ws.Range("B19:C19").Value = 'Duas provas escritas: P1 e P2'
$ws.Range("B20:C20").Value = 'Média das notas obtidas nas duas provas: N1=(P1 + P2)/2'
$ws.Range("B21:C21").Value = 'Uma prova escrita: REC
Média das notas N1 e REC:N2=(N1+REC)/2'
$ws.Range("B22:C22").Value = '1) SMITH, c. A.; CORRIPIO, A. B. Princípios e Prática do Controle Automático de Processo.
3ª ed. Rio de Janeiro: LTC, 2008.
2) STHEPANOPOULOS, G. Chemical Process Control: An Introduction to Theory and Practice. Englewood Cliffs, N.J.: Prentice Hall, 1984.
3) KWONG, W. H. Introdução ao Controle de Processos Químicos com MATLAB. Vols. 1 e 2. São Carlos: EdUFScar, 2002.
4) OGATA, K. Engenharia de Controle Moderno. 5ª ed. São Paulo: Pearson Prentice Hall, 2011.
5) Seborg, D. E.; EDGAR, T. F.; MELLICHAMP, D. A. Process Dynamics and Control. 2 ed. New York: John Wiley & Sons, 2003.
6) COUGHANOWR, D. R.; KOPPEL, L. B. Análise e Controle de Processos. Rio de Janeiro:  Guanabara Dois, 1978.
7) PERRY, R. H.; CHILTON, C. H. Manual de Engenharia Química. 5ª ed. Rio de Janeiro: Guanabara Dois, 1980. Seção 22.'
